# spellcheck & add files
# Fix "radius" -> "diameter" typo in the note column (E2:E4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "could be calculated with 0.95 in diameter, but used 2cm diameter "
$ws.Range("E4").Value = "could be calculated with 0.95 in diameter, but used 2cm diameter "
$ws.Range("E2").Value = "could be calculated with 0.95 in diameter, but used 2cm diameter"

$ws.Range("E3").Select()
